# Apply the documented change to the Blood Suckers 2 review document:
#   1. Insert a new "Meta description" paragraph right after the title
#      (Heading1) paragraph.
#   2. Remove the duplicated bold title paragraph that used to sit near
#      the end of the document (right before the italic meta-description
#      paragraph).
#   3. Replace the text of the remaining (italic) paragraph with the new
#      DALLE feature-image prompt, keeping its italic formatting intact.

$d = $word.ActiveDocument

# Locate (by content, not position) the duplicated bold title paragraph
# and the italic meta-description paragraph near the end of the doc.
$boldTitleIdx = -1
$italicIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Play Blood Suckers 2 Free*" -and $i -ne 1) {
        $boldTitleIdx = $i
    }
    if ($t -like "Read our review of Blood Suckers 2*") {
        $italicIdx = $i
    }
}

# --- 1. Insert "Meta description" paragraph after the title -----------
$title = $d.Paragraphs.Item(1)
$title.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Borrow the "<empty run/><bold run/>" run layout from the bold title
# paragraph further down (same shape we need here), so the new
# paragraph gets a genuine leading empty run like its siblings do.
$boldSource = $d.Paragraphs.Item($boldTitleIdx + 1)
$metaPara.Range.FormattedText = $boldSource.Range.FormattedText

# Re-label the bold run's text, then append the rest of the meta
# description as a plain (non-bold) trailing run.
$boldRun = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$boldRun.Text = "Meta description"
$afterBold = $d.Range($metaPara.Range.End - 1, $metaPara.Range.End - 1)
$afterBold.InsertAfter(": Read our review of Blood Suckers 2, an improved vampire-themed slot machine game from NetEnt. Play for free and enjoy classic Free Spins and a thrilling bonus game.")

# --- 2. Remove the redundant bold title paragraph near the end --------
# (indices shifted down by one now that a paragraph was inserted at the top)
$dupTitle = $d.Paragraphs.Item($boldTitleIdx + 1)
$dupTitle.Range.Delete()

# --- 3. Replace the italic paragraph's text with the DALLE prompt -----
# (index shifted +1 for the inserted paragraph, then -1 back down again
# because the bold-title paragraph right before it was just deleted)
$italicPara = $d.Paragraphs.Item($italicIdx)
$italicRun = $d.Range($italicPara.Range.Start, $italicPara.Range.End - 1)
$italicRun.Text = 'Create a Feature Image Prompt for DALLE: Design a vibrant and cartoonish image featuring a happy and confident Maya warrior, wearing glasses. The warrior should appear to be holding a crossbow while standing confidently in front of a spooky background featuring luxurious coffins. Make sure to include some blood-red elements to emphasize the vampire theme of the game "Blood Suckers 2".'
